# Updated backend and frontend functionality with recent changes
$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Workbook-level: drop workbook protection element
# ------------------------------------------------------------------
$wb.Unprotect() | Out-Null

# ------------------------------------------------------------------
# "Requests" sheet: rewrite rows 2-4, drop rows 5-7
# ------------------------------------------------------------------
$wsReq = $wb.Worksheets.Item("Requests")

# Row 2: ds | 5675 | 2024-09-10 | 2024-09-12 | Yes
$wsReq.Range("A2").Value = "ds"
$wsReq.Range("B2:D2").NumberFormat = "@"
$wsReq.Range("B2").Value = "5675"
$wsReq.Range("C2").Value = "2024-09-10"
$wsReq.Range("D2").Value = "2024-09-12"
$wsReq.Range("E2").Value = "Yes"

# Row 3: ds | 34243 | 2024-09-20 | 2024-09-27 | Yes
$wsReq.Range("A3").Value = "ds"
$wsReq.Range("B3:D3").NumberFormat = "@"
$wsReq.Range("B3").Value = "34243"
$wsReq.Range("C3").Value = "2024-09-20"
$wsReq.Range("D3").Value = "2024-09-27"
$wsReq.Range("E3").Value = "Yes"

# Row 4: SHREE | 677 | 2024-09-10 | 2024-09-19 | Yes
$wsReq.Range("A4").Value = "SHREE"
$wsReq.Range("B4:D4").NumberFormat = "@"
$wsReq.Range("B4").Value = "677"
$wsReq.Range("C4").Value = "2024-09-10"
$wsReq.Range("D4").Value = "2024-09-19"
$wsReq.Range("E4").Value = "Yes"

# Old rows 5-7 no longer exist in the updated sheet
$wsReq.Rows("5:7").Delete() | Out-Null

# ------------------------------------------------------------------
# "Billing" sheet: append rows 2-3
# ------------------------------------------------------------------
$wsBill = $wb.Worksheets.Item("Billing")

# Row 2: hgd | 3244423 | 2024-09-09 | 2024-09-26
$wsBill.Range("A2").Value = "hgd"
$wsBill.Range("B2:D2").NumberFormat = "@"
$wsBill.Range("B2").Value = "3244423"
$wsBill.Range("C2").Value = "2024-09-09"
$wsBill.Range("D2").Value = "2024-09-26"

# Row 3: ds | 34243 | 2024-09-20 | 2024-09-27
$wsBill.Range("A3").Value = "ds"
$wsBill.Range("B3:D3").NumberFormat = "@"
$wsBill.Range("B3").Value = "34243"
$wsBill.Range("C3").Value = "2024-09-20"
$wsBill.Range("D3").Value = "2024-09-27"

# ------------------------------------------------------------------
# View state: restore per-sheet selections and keep "Requests" active
# ------------------------------------------------------------------
$wsBill.Activate() | Out-Null
$wsBill.Range("G19").Select() | Out-Null

$wsReq.Activate() | Out-Null
$wsReq.Range("E31").Select() | Out-Null
